$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge "...войти ил" + [bookmark] + "и зарегистрироваться..." into
#    one continuous run ("...войти или зарегистрироваться...") and
#    fix the "авторизировался" typo + drop "ИГРОВОЙ МАГАЗИН, " from
#    the button list. The Find/Execute below matches across the
#    run/bookmark boundary (Word matches on the logical text stream),
#    and replacing it removes the old, mis-placed "_GoBack" bookmark
#    along with the run split.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "войти ил" + "и зарегистрироваться",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "войти или зарегистрироваться", 2) | Out-Null

# One combined Find/Replace spans from before the misspelled word
# through to after it (past its spellStart/spellEnd proofErr markers)
# and removes "ИГРОВОЙ МАГАЗИН, " from the button list in the same
# pass, so Word collapses everything -- including the orphaned
# proofErr tags -- back into a single plain run, matching the target.
$d.Content.Find.Execute(
    "Если Игрок успешно авторизировался, он попадает в главное меню, " +
    "в котором есть 4 кнопки " + [char]8211 + " ИГРАТЬ, ИГРОВОЙ МАГАЗИН, РЕДАКТОР ГЕРОЯ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Если Игрок успешно авторизовался, он попадает в главное меню, " +
    "в котором есть 4 кнопки " + [char]8211 + " ИГРАТЬ, РЕДАКТОР ГЕРОЯ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Delete the whole "ИГРОВОЙ МАГАЗИН" paragraph (now paragraph 3,
#    1-based) -- the in-game shop feature was cut from the spec.
# ------------------------------------------------------------------
$shopPara = $d.Paragraphs.Item(3)
$shopPara.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the end of the "ВЫХОД"
#    paragraph (paragraph 2), right after "закрывается. " and before
#    the paragraph mark. A collapsed range placed directly at that
#    offset is mis-resolved by Bookmarks.Add, so we work around it:
#    insert a 1-char placeholder there, wrap a (non-collapsed)
#    bookmark range around the placeholder, then erase the
#    placeholder through the bookmark's own range so the bookmark
#    collapses back down to the correct spot.
# ------------------------------------------------------------------
$exitPara = $d.Paragraphs.Item(2)
$exitEnd = $exitPara.Range.End - 1

$placeholder = $d.Range($exitEnd, $exitEnd)
$placeholder.InsertAfter([char]8203)

$wrap = $d.Range($exitEnd, $exitEnd + 1)
$d.Bookmarks.Add("_GoBack", $wrap) | Out-Null
$goBack = $d.Bookmarks("_GoBack")
$goBack.Range.Text = ""

# ------------------------------------------------------------------
# 4) Drop the in-game-currency reward clause from the closing
#    paragraph -- the shop is gone, so the win reward is just XP now.
#    The whole paragraph text is matched/replaced (rather than just
#    the removed clause) so Word collapses it back into a single run,
#    clearing the mid-paragraph lastRenderedPageBreak + proofErr
#    markers that littered the original multi-run paragraph.
# ------------------------------------------------------------------
$dash = [char]8211
$lastOld = "Игра заключается в возможности ходить по квадратам не являющимися стенами и расставлять бомбы, которые взрываются по прошествии определенного времени. В зависимости от выбранных игроком навыков, бомбы могут взрываться по-разному, уничтожая находящиеся вокруг них стены и повреждая игроков. Победителем считается последний выживший игрок. Победителю начисляется определенное количество внутриигровой валюты, которыю можно тратить во внутриигровом магазине, а также некоторое количество опыта. При достижении определенного (для каждого уровня героя своего) количества опыта, происходит повышение уровня героя, у героя появляется новое очко навыков. Если несколько игроков остались живы по истечению определенного времени $dash времени раунда, все игроки считаются проигравшими."
$lastNew = "Игра заключается в возможности ходить по квадратам не являющимися стенами и расставлять бомбы, которые взрываются по прошествии определенного времени. В зависимости от выбранных игроком навыков, бомбы могут взрываться по-разному, уничтожая находящиеся вокруг них стены и повреждая игроков. Победителем считается последний выживший игрок. Победителю начисляется некоторое количество опыта. При достижении определенного (для каждого уровня героя своего) количества опыта, происходит повышение уровня героя, у героя появляется новое очко навыков. Если несколько игроков остались живы по истечению определенного времени $dash времени раунда, все игроки считаются проигравшими."

$d.Content.Find.Execute(
    $lastOld, $true, $false, $false, $false, $false, $true, 1, $false,
    $lastNew, 2) | Out-Null
